$wb = $excel.ActiveWorkbook
$wsView = $wb.Worksheets.Item("view")
$wsCtrl = $wb.Worksheets.Item("controleur")

# --- Content changes on the "controleur" sheet -----------------------------
# Row 5: "livrer_commande" -> "delivreryOrder", mark "Fait" column with V
$wsCtrl.Range("A5").Value = "delivreryOrder"
$wsCtrl.Range("D5").Value = "V"

# Row 6: "ajouter_produit" -> "addOrderDetail", mark "Fait" column with V
$wsCtrl.Range("A6").Value = "addOrderDetail"
$wsCtrl.Range("D6").Value = "V"

# Row 7: mark "Fait" column with V
$wsCtrl.Range("D7").Value = "V"

# Row 8: "supprimer_utilisateur" entry removed entirely (row cleared)
$wsCtrl.Range("A8:C8").ClearContents()

# Row 9: "afficher_detail_commande" -> "displayDetailOrder", new parameter text, mark V
$wsCtrl.Range("A9").Value = "displayDetailOrder"
$wsCtrl.Range("C9").Value = "`$listOrder(liste insexé par l'id), `$detailOrder(liste Indexé par l'id); `$orderCurrent(comande selectionné)"
$wsCtrl.Range("D9").Value = "V"

# Row 12: mark "Fait" column with V
$wsCtrl.Range("D12").Value = "V"

# Row 18: mark "Fait" column with V
$wsCtrl.Range("D18").Value = "V"

# --- Sheet view / selection state -------------------------------------------
# controleur: selection moves from A2 to D18, and the stray topLeftCell scroll
# position is dropped when we select on this sheet.
$wsCtrl.Range("D18").Select()

# view: becomes the active / selected tab, keeping its existing C7 selection.
$wsView.Activate()
$wsView.Range("C7").Select()
